# Add mock data for last month
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 0.75
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("AG3").Value = 0.8
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 0.8
$ws.Range("AH4").Value = 0.7
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("AG5").Value = 0.6
$ws.Range("AI5").ClearContents()

# Header row - clear the date value in AI1 but keep its style/formatting
$ws.Range("AI1").ClearContents()

# Update the active selection to match the saved view state
$ws.Range("AG7").Select()
